# Daily attendance processing - 2026-01-14 07:16:44
# Swap the order of "System" and the recorder's email address in the
# "Recorded By" column (G) wherever the cell currently reads
# "System, dnasr281@gmail.com", changing it to
# "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

# Used range dimensions so we scan every data row under the "Recorded By"
# header in column G (col 7), rather than hard-coding row numbers.
$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

$updated = 0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
        $updated = $updated + 1
    }
}

Write-Host "Updated $updated 'Recorded By' cell(s) in column G"
